$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '98.804.38'
$ws.Cells.Item(2, 5).Value = '  +1.58%  '

$ws.Cells.Item(3, 4).Value = '3.306.05'
$ws.Cells.Item(3, 5).Value = '  -0.88%  '

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '255.33'
$ws.Cells.Item(5, 5).Value = '  +0.46%  '

$ws.Cells.Item(6, 4).Value = '625.04'
$ws.Cells.Item(6, 5).Value = '  +0.73%  '

$ws.Cells.Item(7, 5).Value = '  +33.19%  '

$ws.Cells.Item(8, 4).Value = '0.403'
$ws.Cells.Item(8, 5).Value = '  +5.06%  '

$ws.Cells.Item(9, 4).Value = '0.998'
$ws.Cells.Item(9, 5).Value = '  -0.09%  '

$ws.Cells.Item(10, 4).Value = '0.972'
$ws.Cells.Item(10, 5).Value = '  +23.81%  '

$ws.Cells.Item(11, 4).Value = '3.303.85'
$ws.Cells.Item(11, 5).Value = '  -0.94%  '

$ws.Cells.Item(13, 4).Value = '39.73'
$ws.Cells.Item(13, 5).Value = '  +12.50%  '

$ws.Cells.Item(14, 4).Value = '98.472.59'
$ws.Cells.Item(14, 5).Value = '  +1.42%  '

$ws.Cells.Item(15, 4).Value = '0.0000249'
$ws.Cells.Item(15, 5).Value = '  +1.39%  '

$ws.Cells.Item(16, 4).Value = '3.914.65'
$ws.Cells.Item(16, 5).Value = '  -0.05%  '

$ws.Cells.Item(17, 5).Value = '  -0.86%  '

$ws.Cells.Item(18, 4).Value = '3.310.69'
$ws.Cells.Item(18, 5).Value = '  -0.56%  '

$ws.Cells.Item(19, 5).Value = '  -2.66%  '

$ws.Cells.Item(20, 4).Value = '15.63'
$ws.Cells.Item(20, 5).Value = '  +4.89%  '

$ws.Cells.Item(21, 4).Value = '6.31'
$ws.Cells.Item(21, 5).Value = '  +9.04%  '

$ws.Cells.Item(22, 4).Value = '487.19'
$ws.Cells.Item(22, 5).Value = '  +1.20%  '

$ws.Cells.Item(23, 4).Value = '9.47'
$ws.Cells.Item(23, 5).Value = '  +2.74%  '

$ws.Cells.Item(24, 5).Value = '  -3.63%  '

$ws.Cells.Item(25, 4).Value = '5.61'
$ws.Cells.Item(25, 5).Value = '  -0.70%  '

$ws.Cells.Item(26, 4).Value = '88.56'
$ws.Cells.Item(26, 5).Value = '  +1.09%  '

$ws.Cells.Item(27, 4).Value = '12.03'
$ws.Cells.Item(27, 5).Value = '  -0.02%  '

$ws.Cells.Item(28, 5).Value = '  +27.07%  '

$ws.Cells.Item(29, 4).Value = '3.477.49'
$ws.Cells.Item(29, 5).Value = '  -0.75%  '

$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  -0.11%  '

$ws.Cells.Item(31, 4).Value = '0.135'
$ws.Cells.Item(31, 5).Value = '  +11.58%  '

$ws.Cells.Item(32, 5).Value = '  +3.36%  '

$ws.Cells.Item(33, 4).Value = '10.05'
$ws.Cells.Item(33, 5).Value = '  +9.40%  '

$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 5).Value = '  +0.16%  '

$ws.Cells.Item(35, 4).Value = '28.06'
$ws.Cells.Item(35, 5).Value = '  +3.00%  '

$ws.Cells.Item(36, 2).Value = 'RenderToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(36, 4).Value = '7.18'
$ws.Cells.Item(36, 5).Value = '  -2.85%  '

$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).Value = '0.147'
$ws.Cells.Item(37, 5).Value = '  -2.19%  '

$ws.Cells.Item(38, 4).Value = '1.95'
$ws.Cells.Item(38, 5).Value = '  +0.73%  '

$ws.Cells.Item(39, 4).Value = '0.466'
$ws.Cells.Item(39, 5).Value = '  +3.87%  '

$ws.Cells.Item(40, 4).Value = '24.85'
$ws.Cells.Item(40, 5).Value = '  +0.29%  '

$ws.Cells.Item(41, 4).Value = '489.34'
$ws.Cells.Item(41, 5).Value = '  -3.70%  '

$ws.Cells.Item(42, 4).Value = '3.66'
$ws.Cells.Item(42, 5).Value = '  +4.01%  '

$ws.Cells.Item(43, 5).Value = '  -3.36%  '

$ws.Cells.Item(44, 4).Value = '0.788'
$ws.Cells.Item(44, 5).Value = '  -0.68%  '

$ws.Cells.Item(45, 5).Value = '  -0.01%  '

$ws.Cells.Item(46, 4).Value = '3.11'
$ws.Cells.Item(46, 5).Value = '  -5.72%  '

$ws.Cells.Item(47, 5).Value = '  -1.28%  '

$ws.Cells.Item(48, 4).Value = '1.94'
$ws.Cells.Item(48, 5).Value = '  +1.44%  '

$ws.Cells.Item(49, 4).Value = '7.33'
$ws.Cells.Item(49, 5).Value = '  +16.36%  '

$ws.Cells.Item(50, 4).Value = '0.847'
$ws.Cells.Item(50, 5).Value = '  +6.00%  '

$ws.Cells.Item(51, 5).Value = '  +5.10%  '
